# AT: script changes to omit "Other" land from statistics
#
# The "Otherland_pcnt" row (row 5) is dropped entirely from the county-means
# summary table (its underlying data is omitted upstream), which shifts all
# subsequent rows up by one. Because the land-use percentages are shares of
# a total that no longer includes "Other" land, the percentages for the
# remaining land-use categories (Cropland, CRPland, Forestland, Pastureland,
# Rangeland, Urbanland) are recomputed with the new (smaller) denominator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Otherland_pcnt" row (row 5); everything below shifts up.
$ws.Rows.Item(5).Delete()

# Recomputed percentages (new denominator excludes "Other" land) for the
# rows that used to be 2-7 (Cropland_pcnt, CRPland_pcnt, Forestland_pcnt,
# Pastureland_pcnt, Rangeland_pcnt, Urbanland_pcnt) and now occupy rows 2-7
# after the shift.

$newValues = @{
    2 = @(28.125761032104492, 26.64579963684082, 25.165998458862305, 24.52128791809082, 23.52459716796875, 23.244480133056641, 23.984085083007813)
    3 = @(0, 0.85870075225830078, 2.1823241710662842, 2.0761411190032959, 1.8803716897964478, 1.9339430332183838, 1.4021095037460327)
    4 = @(34.908100128173828, 35.098213195800781, 35.120071411132813, 35.134029388427734, 35.237506866455078, 35.257175445556641, 35.342998504638672)
    5 = @(16.106454849243164, 16.148845672607422, 15.718827247619629, 15.413778305053711, 15.694674491882324, 15.392695426940918, 14.879504203796387)
    6 = @(15.54348087310791, 15.309108734130859, 15.159688949584961, 15.129592895507813, 15.120001792907715, 15.093488693237305, 15.042671203613281)
    7 = @(5.3162059783935547, 5.9393305778503418, 6.6530900001525879, 7.7251725196838379, 8.5428485870361328, 9.0782184600830078, 9.3486337661743164)
}

$columns = @(2, 3, 4, 5, 6, 7, 8)  # B..H

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Cells.Item($row, $columns[$i]).Value = $vals[$i]
    }
}
